# Removed unused master slide
#
# The "Title Slide" custom layout (ppt/slideLayouts/slideLayout18.xml,
# referenced as rId18 / sldLayoutId 2147483723 from slideMaster1.xml) is not
# used by any slide in the deck, so it is deleted here. Deleting it removes
# the layout part, its relationship, the Content_Types override, and the
# <p:sldLayoutId> entry that pointed at it.

$p = $ppt.ActivePresentation

# Resolve the slide master through Designs (Design.SlideMaster) so that
# CustomLayouts.Item(n) addresses layouts by their real, stable index.
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

# Find the "Title Slide" layout (the one being removed) by name, rather than
# assuming it is always the last one.
$layoutToRemove = $null
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $candidate = $master.CustomLayouts.Item($i)
    if ($candidate.Name -eq "Title Slide") {
        $layoutToRemove = $candidate
    }
}

if ($layoutToRemove -ne $null) {
    $layoutToRemove.Delete()
}
